$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("2018-10", 101.2, 102.2, 101),
  @("2018-11", 101, 101.9, 101),
  @("2018-12", 101.3, 101.6, 101.2),
  @("2018-01", 101.3532, 101.9457, 100.714),
  @("2018-02", 101.7, 102, 100.9),
  @("2018-03", 101.8, 102.2, 101),
  @("2018-04", 101.4, 102.2, 101),
  @("2018-05", 101.1, 102.1, 100.9),
  @("2018-06", 100.9, 102.6, 100.9),
  @("2018-07", 101.9, 102.5, 100.7),
  @("2018-08", 101.5, 102.2, 100.8),
  @("2018-09", 101.4, 102.2, 100.8),
  @("2019-10", 101.8, 101.2, 101.6),
  @("2019-11", 102.1, 101, 101.5),
  @("2019-12", 101.5, 101.1, 100.9),
  @("2019-01", 100.9, 101, 100.9),
  @("2019-02", 100.4, 100.9, 100.8),
  @("2019-03", 100.5, 100.8, 100.8),
  @("2019-04", 101.6, 100.8, 101.2),
  @("2019-05", 101.1, 101, 101.3),
  @("2019-06", 101.2, 101.3, 101.4),
  @("2019-07", 101.2, 101.4, 101.7),
  @("2019-08", 101.8, 101.5, 101.8),
  @("2019-09", 101.6, 101.3, 101.8),
  @("2020-10", 99.6, 101.5, 98.8),
  @("2020-11", 99.3, 102.2, 98.8),
  @("2020-12", 99.3, 102.3, 99),
  @("2020-01", 101.5, 101.3, 100.7),
  @("2020-02", 101.3, 101.5, 100.8),
  @("2020-03", 101, 101.5, 100.1),
  @("2020-04", 99.4, 101.4, 99.9),
  @("2020-05", 99.2, 101.4, 99.7),
  @("2020-06", 99.1, 101.1, 99.7),
  @("2020-07", 99.4, 101.3, 99.5),
  @("2020-08", 99.5, 101.4, 99),
  @("2020-09", 99.5, 101.6, 98.9),
  @("2021-10", 99.4, 102.4, 100),
  @("2021-11", 99.2, 102.1, 99.9),
  @("2021-12", 99.7, 101.5, 100.1),
  @("2021-01", 99.7, 103.3, 99.5),
  @("2021-02", 99.9, 103.5, 99.6),
  @("2021-03", 100.1, 103.7, 100),
  @("2021-04", 100.3, 103.8, 99.6),
  @("2021-05", 100.3, 103.7, 99.7),
  @("2021-06", 100.3, 103.5, 99.9),
  @("2021-07", 100.2, 103.4, 99.8),
  @("2021-08", 100, 102.8, 99.7),
  @("2021-09", 100.1, 102.8, 99.6),
  @("2022-10", 100.8, 101.2, 102.3),
  @("2022-11", 101.3, 101, 102.1),
  @("2022-12", 101.1, 100.9, 102.1),
  @("2022-01", 99.2, 101.3, 100.5),
  @("2022-02", 99.2, 100.6, 100.8),
  @("2022-03", 99.1, 101, 100.8),
  @("2022-04", 100.4, 101, 101.1),
  @("2022-05", 100.4, 100.8, 101.3),
  @("2022-06", 100.7, 100.7, 101.3),
  @("2022-07", 99.9, 100.5, 101.1),
  @("2022-08", 100.1, 101.1, 101.8),
  @("2022-09", 100.2, 101, 102.2),
  @("2023-01", 101, 101.3, 102.1),
  @("2023-02", 101.3, 101.3, 102.5),
  @("2023-03", 102.1, 100.4, 102.7),
  @("2023-04", 101.4, 99.9, 102.6),
  @("2023-05", 101.5, 100.5, 101.9),
  @("2023-06", 101.4, 100.6, 101.9),
  @("2023-07", 102, 100.9, 102)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Apply date-column style (border/bold/center) to newly added rows 50-68
$ws.Range("A2").Copy()
$ws.Range("A50:A68").PasteSpecial(-4122)